$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price and Volume columns retain their original text formatting
# (values like "1.009", "2.330", "0.00001044" must not be reinterpreted as numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '21.067.24'
$ws.Range('E2').Value = '  +3.23%  '

$ws.Range('D3').Value = '1.534.35'
$ws.Range('E3').Value = '  +4.93%  '

$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  +0.45%  '

$ws.Range('D5').Value = '0.9628'
$ws.Range('E5').Value = '  +1.28%  '

$ws.Range('D6').Value = '281.85'
$ws.Range('E6').Value = '  +2.59%  '

$ws.Range('D7').Value = '0.3622'
$ws.Range('E7').Value = '  -0.82%  '

$ws.Range('D8').Value = '0.3187'
$ws.Range('E8').Value = '  +3.81%  '

$ws.Range('D9').Value = '40.65'
$ws.Range('E9').Value = '  +2.06%  '

$ws.Range('D10').Value = '1.091'
$ws.Range('E10').Value = '  +5.51%  '

$ws.Range('D11').Value = '0.06794'
$ws.Range('E11').Value = '  +3.26%  '

$ws.Range('D12').Value = '1.008'
$ws.Range('E12').Value = '  +0.69%  '

$ws.Range('D13').Value = '5.666'
$ws.Range('E13').Value = '  +4.41%  '

$ws.Range('D14').Value = '18.67'
$ws.Range('E14').Value = '  +4.00%  '

$ws.Range('D15').Value = '6.346'
$ws.Range('E15').Value = '  +3.36%  '

$ws.Range('D16').Value = '0.00001044'
$ws.Range('E16').Value = '  +2.03%  '

$ws.Range('D17').Value = '0.9629'
$ws.Range('E17').Value = '  -0.79%  '

$ws.Range('D18').Value = '1.525.96'
$ws.Range('E18').Value = '  +4.39%  '

$ws.Range('D19').Value = '0.06086'
$ws.Range('E19').Value = '  +4.66%  '

$ws.Range('D20').Value = '72.19'
$ws.Range('E20').Value = '  +4.08%  '

$ws.Range('D21').Value = '5.696'
$ws.Range('E21').Value = '  +4.78%  '

$ws.Range('D22').Value = '14.96'
$ws.Range('E22').Value = '  +3.63%  '

$ws.Range('D23').Value = '11.34'
$ws.Range('E23').Value = '  +4.06%  '

$ws.Range('D24').Value = '2.330'
$ws.Range('E24').Value = '  +3.72%  '

$ws.Range('D25').Value = '21.178.32'
$ws.Range('E25').Value = '  +3.63%  '

$ws.Range('D26').Value = '148.26'
$ws.Range('E26').Value = '  +4.65%  '

$ws.Range('D27').Value = '2.208'
$ws.Range('E27').Value = '  +5.98%  '

$ws.Range('D28').Value = '17.69'
$ws.Range('E28').Value = '  +3.39%  '

$ws.Range('D29').Value = '1.696.83'
$ws.Range('E29').Value = '  +5.14%  '

$ws.Range('D30').Value = '118.23'
$ws.Range('E30').Value = '  +5.27%  '

$ws.Range('D31').Value = '4.032'
$ws.Range('E31').Value = '  +4.36%  '

$ws.Range('D32').Value = '0.8512'
$ws.Range('E32').Value = '  +7.79%  '

$ws.Range('D33').Value = '5.179'
$ws.Range('E33').Value = '  +5.94%  '

$ws.Range('D34').Value = '0.07998'
$ws.Range('E34').Value = '  +1.27%  '

$ws.Range('E35').Value = '  -1.41%  '

$ws.Range('D36').Value = '4.957'
$ws.Range('E36').Value = '  +5.89%  '

$ws.Range('D37').Value = '1.206'
$ws.Range('E37').Value = '  +4.88%  '

$ws.Range('D38').Value = '0.05857'
$ws.Range('E38').Value = '  +2.30%  '

$ws.Range('D39').Value = '0.02102'
$ws.Range('E39').Value = '  +3.78%  '

$ws.Range('D40').Value = '10.65'
$ws.Range('E40').Value = '  +3.18%  '

$ws.Range('D41').Value = '7.709'
$ws.Range('E41').Value = '  +3.35%  '

$ws.Range('D42').Value = '0.1914'
$ws.Range('E42').Value = '  +3.02%  '

$ws.Range('D43').Value = '0.9632'
$ws.Range('E43').Value = '  +0.54%  '

$ws.Range('D44').Value = '0.5427'
$ws.Range('E44').Value = '  +3.19%  '

$ws.Range('E45').Value = '  +4.93%  '

$ws.Range('D46').Value = '3.577'
$ws.Range('E46').Value = '  +2.55%  '

$ws.Range('D47').Value = '0.5446'
$ws.Range('E47').Value = '  +6.02%  '

$ws.Range('D48').Value = '121.24'
$ws.Range('E48').Value = '  +3.58%  '

$ws.Range('D49').Value = '1.865'
$ws.Range('E49').Value = '  +6.78%  '

$ws.Range('D50').Value = '0.06566'
$ws.Range('E50').Value = '  +2.29%  '

$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = '0.9926'
$ws.Range('E51').Value = '  +0.08%  '
